$d = $word.ActiveDocument

$old = "workforce on top of increasing prices for goods and services, making it"
$new = "workforce and increase prices" + [char]0x2014 + "making it"

$range = $d.Content
$range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
